# Update "想去人数" (want-to-go count) figures for the 展览 (Exhibition) and
# 全部类型 (All Types) sheets, reflecting the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row number -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6363
$ws1.Range("F6").Value = 55
$ws1.Range("F9").Value = 63
$ws1.Range("F12").Value = 149
$ws1.Range("F13").Value = 358
$ws1.Range("F14").Value = 619
$ws1.Range("F15").Value = 3098
$ws1.Range("F18").Value = 1755

# Sheet "全部类型": row number -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6363
$ws4.Range("F6").Value = 55
$ws4.Range("F10").Value = 63
$ws4.Range("F13").Value = 149
$ws4.Range("F14").Value = 358
$ws4.Range("F15").Value = 619
$ws4.Range("F16").Value = 3098
$ws4.Range("F19").Value = 1755
